$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the existing "LUIS DANIEL DE LA ROSA PEREZ"
# row (row 16), pushing the old row 17 down to row 19 and copying row 16's
# formatting (borders/number formats) into the two freshly inserted rows so
# they render like the rest of the data table rather than with blank formats.
$ws.Rows("17:18").Insert()
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# Row 16: LUIS DANIEL DE LA ROSA PEREZ, periodo 2507, nueva mora 56940
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1007126667"
$ws.Range("D16").Value = "LUIS DANIEL DE LA ROSA PEREZ"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# Row 17: LUIS DANIEL DE LA ROSA PEREZ, periodo 2505 (pre-existing record)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1007126667"
$ws.Range("D17").Value = "LUIS DANIEL DE LA ROSA PEREZ"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 13286
$ws.Range("G17").Value = 1423500

# Row 18: JESUS ALBERTO ESCORCIA SANTIAGO, periodo 2507, nueva mora 56940
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1235039810"
$ws.Range("D18").Value = "JESUS ALBERTO ESCORCIA SANTIAGO"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19 (previously row 17): JESUS ALBERTO ESCORCIA SANTIAGO, periodo 2505
# -- content unchanged, just shifted down by the two inserted rows above.
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1235039810"
$ws.Range("D19").Value = "JESUS ALBERTO ESCORCIA SANTIAGO"
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 13286
$ws.Range("G19").Value = 1423500

# Cant. Periodos goes from 1 to 2 (two periods per worker now).
$ws.Range("F13").Value = 2

# Valor Mora total is the sum of all "Valor Mora" entries in the table.
$ws.Range("E11").Value = 140452
